$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2466.3333
$ws.Range("I4").Value = 2466.3333
$ws.Range("K4").Value = 2466.3333
$ws.Range("M4").Value = -2352.3333
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("H17").Value = 1780
$ws.Range("I17").Value = 1950
$ws.Range("J17").Value = 1737.5
$ws.Range("K17").Value = 5850
$ws.Range("L17").Value = 5212.5
$ws.Range("M17").Value = -5682
$ws.Range("N17").Value = -5548.5
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H70").Value = 1281.3636
$ws.Range("I70").Value = 865
$ws.Range("J70").Value = 1437.5
$ws.Range("K70").Value = 2595
$ws.Range("L70").Value = 4312.5
$ws.Range("M70").Value = -2325
$ws.Range("N70").Value = -4852.5
$ws.Range("H73").Value = 1281.3636
$ws.Range("I73").Value = 865
$ws.Range("J73").Value = 1437.5
$ws.Range("K73").Value = 2595
$ws.Range("L73").Value = 4312.5
$ws.Range("M73").Value = -1659
$ws.Range("N73").Value = -6184.5
$ws.Range("H74").Value = 4750
$ws.Range("I74").Value = 4666.6665
$ws.Range("K74").Value = 4666.6665
$ws.Range("M74").Value = -3730.6665
$ws.Range("H77").Value = 4750
$ws.Range("I77").Value = 4666.6665
$ws.Range("K77").Value = 23333.3325
$ws.Range("M77").Value = -18653.3325
$ws.Range("H92").Value = 3141
$ws.Range("I92").Value = 2747.5
$ws.Range("J92").Value = 3665.6667
$ws.Range("K92").Value = 2747.5
$ws.Range("L92").Value = 3665.6667
$ws.Range("M92").Value = -1499.5
$ws.Range("N92").Value = -6161.6667
$ws.Range("H137").Value = 2350.4285
$ws.Range("I137").Value = 2546.5
$ws.Range("K137").Value = 7639.5
$ws.Range("M137").Value = -5089.5
$ws.Range("H138").Value = 3430.9524
$ws.Range("I138").Value = 1016.75
$ws.Range("K138").Value = 3050.25
$ws.Range("M138").Value = 2089.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2898.6365
$ws.Range("I61").Value = 2898.6365
$ws.Range("K61").Value = 2898.6365
$ws.Range("M61").Value = -2686.6365
$ws.Range("H74").Value = 2156.182
$ws.Range("I74").Value = 1777.25
$ws.Range("K74").Value = 1777.25
$ws.Range("M74").Value = -903.25
$ws.Range("H77").Value = 2156.182
$ws.Range("I77").Value = 1777.25
$ws.Range("K77").Value = 8886.25
$ws.Range("M77").Value = -4518.25
$ws.Range("H136").Value = 2898.6365
$ws.Range("I136").Value = 2898.6365
$ws.Range("K136").Value = 8695.9095
$ws.Range("M136").Value = -6145.9095

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 671.6667
$ws.Range("I5").Value = 671.6667
$ws.Range("K5").Value = 671.6667
$ws.Range("M5").Value = -558.6667
$ws.Range("H30").Value = 10000
$ws.Range("I30").Value = 10000
$ws.Range("K30").Value = 10000
$ws.Range("M30").Value = -9875

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 600
$ws.Range("I16").Value = 600
$ws.Range("J16").Value = 600
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = -313
$ws.Range("N16").Value = -1174
$ws.Range("H31").Value = 6192.722
$ws.Range("I31").Value = 5891.143
$ws.Range("J31").Value = 7248.25
$ws.Range("K31").Value = 5891.143
$ws.Range("L31").Value = 7248.25
$ws.Range("M31").Value = -5596.143
$ws.Range("N31").Value = -7838.25
$ws.Range("H34").Value = 6192.722
$ws.Range("I34").Value = 5891.143
$ws.Range("J34").Value = 7248.25
$ws.Range("K34").Value = 5891.143
$ws.Range("L34").Value = 7248.25
$ws.Range("M34").Value = -5689.143
$ws.Range("N34").Value = -7652.25
$ws.Range("H47").Value = 9999
$ws.Range("J47").Value = 9999
$ws.Range("L47").Value = 9999
$ws.Range("N47").Value = -11131
$ws.Range("H58").Value = 2725.6667
$ws.Range("I58").Value = 2725.6667
$ws.Range("K58").Value = 2725.6667
$ws.Range("M58").Value = -2522.6667
$ws.Range("H105").Value = 1499
$ws.Range("I105").Value = 1499
$ws.Range("K105").Value = 1499
$ws.Range("M105").Value = 248
$ws.Range("H107").Value = 9000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws.Range("H113").Value = 600
$ws.Range("I113").Value = 600
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 600
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = 1570
$ws.Range("N113").Value = -4940
$ws.Range("H132").Value = 1965.1666
$ws.Range("I132").Value = 1965.1666
$ws.Range("K132").Value = 5895.4998
$ws.Range("M132").Value = -3365.4998
$ws.Range("H136").Value = 2725.6667
$ws.Range("I136").Value = 2725.6667
$ws.Range("K136").Value = 8177.000100000001
$ws.Range("M136").Value = -5627.000100000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 32.882355
$ws.Range("I2").Value = 5.909091
$ws.Range("J2").Value = 82.333336
$ws.Range("K2").Value = 35.454546
$ws.Range("L2").Value = 494.000016
$ws.Range("M2").Value = 77.54545400000001
$ws.Range("N2").Value = -720.000016
$ws.Range("H10").Value = 2531.375
$ws.Range("I10").Value = 35.857143
$ws.Range("J10").Value = 20000
$ws.Range("K10").Value = 107.571429
$ws.Range("L10").Value = 60000
$ws.Range("M10").Value = 31.42857100000001
$ws.Range("N10").Value = -60278
$ws.Range("H26").Value = 60
$ws.Range("I26").Value = 100
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = 300
$ws.Range("L26").Value = 60
$ws.Range("M26").Value = -12
$ws.Range("N26").Value = -636
$ws.Range("H103").Value = 2845
$ws.Range("I103").Value = 2560
$ws.Range("J103").Value = 3082.5
$ws.Range("K103").Value = 7680
$ws.Range("L103").Value = 9247.5
$ws.Range("M103").Value = -6801
$ws.Range("N103").Value = -11005.5
$ws.Range("H131").Value = 4119.231
$ws.Range("I131").Value = 2433.3333
$ws.Range("J131").Value = 4625
$ws.Range("K131").Value = 7299.999899999999
$ws.Range("L131").Value = 13875
$ws.Range("M131").Value = -2259.999899999999
$ws.Range("N131").Value = -23955
$ws.Range("H132").Value = 989.3333
$ws.Range("I132").Value = 989.3333
$ws.Range("K132").Value = 8903.9997
$ws.Range("M132").Value = -6373.9997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 834500
$ws.Range("I3").Value = 834500
$ws.Range("K3").Value = 834500
$ws.Range("M3").Value = -834384
$ws.Range("H40").Value = 30000
$ws.Range("I40").Value = 30000
$ws.Range("K40").Value = 30000
$ws.Range("M40").Value = -29849
$ws.Range("H43").Value = 8550
$ws.Range("I43").Value = 6937.5
$ws.Range("J43").Value = 15000
$ws.Range("K43").Value = 6937.5
$ws.Range("L43").Value = 15000
$ws.Range("M43").Value = -6786.5
$ws.Range("N43").Value = -15302
$ws.Range("H102").Value = 2690.7
$ws.Range("I102").Value = 2189
$ws.Range("K102").Value = 2189
$ws.Range("M102").Value = -567
$ws.Range("H126").Value = 7949.3335
$ws.Range("I126").Value = 5924
$ws.Range("J126").Value = 12000
$ws.Range("K126").Value = 17772
$ws.Range("L126").Value = 36000
$ws.Range("M126").Value = -15302
$ws.Range("N126").Value = -40940
$ws.Range("H132").Value = 3132.5386
$ws.Range("I132").Value = 2099.8572
$ws.Range("K132").Value = 6299.571599999999
$ws.Range("M132").Value = -3769.571599999999
$ws.Range("H138").Value = 200000
$ws.Range("I138").Value = 200000
$ws.Range("K138").Value = 200000
$ws.Range("M138").Value = -194860

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 676.5
$ws.Range("I9").Value = 569
$ws.Range("K9").Value = 569
$ws.Range("M9").Value = -345
$ws.Range("H55").Value = 1459.8
$ws.Range("I55").Value = 324.75
$ws.Range("J55").Value = 6000
$ws.Range("K55").Value = 324.75
$ws.Range("L55").Value = 6000
$ws.Range("M55").Value = -151.75
$ws.Range("N55").Value = -6346
$ws.Range("H61").Value = 2292.2
$ws.Range("I61").Value = 2292.2
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2292.2
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2090.2
$ws.Range("N61").ClearContents()
$ws.Range("H113").Value = 2292.2
$ws.Range("I113").Value = 2292.2
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2292.2
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -122.1999999999998
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 3812.6667
$ws.Range("J132").Value = 1250
$ws.Range("L132").Value = 3750
$ws.Range("N132").Value = -8810

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 40000000
$ws.Range("I10").Value = 40000000
$ws.Range("K10").Value = 40000000
$ws.Range("M10").Value = -39999831
$ws.Range("H136").Value = 1021.5714
$ws.Range("I136").Value = 1021.5714
$ws.Range("K136").Value = 3064.7142
$ws.Range("M136").Value = -514.7142000000003

Write-Host "Applied all Kraken_Profits updates"